# Update the "想去人数" (number of people interested) figures in both the
# "展览" and "全部类型" sheets to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 2971
    7  = 1619
    10 = 27
    11 = 1328
    13 = 448
    15 = 69
    19 = 100
    20 = 3051
    22 = 88
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
